$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.091.48"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.597.72"
$ws.Range("E3").Value = "  +1.97%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.87"
$ws.Range("E5").Value = "  +2.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.06"
$ws.Range("E6").Value = "  +1.38%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +2.90%  "
$ws.Range("E9").Value = "  +3.11%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.69"
$ws.Range("E10").Value = "  +3.42%  "
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("E12").Value = "  +0.02%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.27"
$ws.Range("E13").Value = "  +0.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.063.65"
$ws.Range("E14").Value = "  +2.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "62.968.31"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("E16").Value = "  +3.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.597.56"
$ws.Range("E17").Value = "  +1.81%  "
$ws.Range("E18").Value = "  +0.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "343.39"
$ws.Range("E19").Value = "  +2.63%  "
$ws.Range("E20").Value = "  +1.42%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.79"
$ws.Range("E21").Value = "  +0.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.00"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.73"
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.22"
$ws.Range("E24").Value = "  +2.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.726.41"
$ws.Range("E25").Value = "  +2.00%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("E27").Value = "  +1.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").Value = "  -0.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.41"
$ws.Range("E29").Value = "  +0.88%  "
$ws.Range("E30").Value = "  +8.02%  "
$ws.Range("E31").Value = "  -0.05%  "
$ws.Range("E32").Value = "  +5.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0₃0823"
$ws.Range("E33").Value = "  +1.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "466.05"
$ws.Range("E34").Value = "  +15.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "176.84"
$ws.Range("E35").Value = "  +0.72%  "
$ws.Range("E36").Value = "  +4.43%  "
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("E38").Value = "  +0.67%  "
$ws.Range("E39").Value = "  +0.84%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.60"
$ws.Range("E40").Value = "  +6.05%  "
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("E42").Value = "  -1.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "159.02"
$ws.Range("E43").Value = "  +4.84%  "
$ws.Range("E44").Value = "  +1.54%  "
$ws.Range("E45").Value = "  +6.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "21.14"
$ws.Range("E46").Value = "  +2.09%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0548"
$ws.Range("E47").Value = "  +3.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0975"
$ws.Range("E48").Value = "  +1.09%  "
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.58"
$ws.Range("E50").Value = "  +1.77%  "
$ws.Range("E51").Value = "  +0.35%  "
